$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 73

$ws.Cells.Item($row, 1).Value = "2024-10-13 00:00:00"
$ws.Cells.Item($row, 2).Value = 76450
$ws.Cells.Item($row, 3).Value = 10776.11
$ws.Cells.Item($row, 4).Value = 9536.38
$ws.Cells.Item($row, 5).Value = 7.0662
